$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / URL / name / percentage cells - direct assignment is safe
# (Excel cannot coerce these strings into numbers because of stray dots,
#  spaces, or non-numeric characters).
$ws.Range("D2").Value = "26.627.94"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.854.03"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "1.851.51"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "26.663.73"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  -7.18%  "
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +4.93%  "
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("E42").Value = "  -3.28%  "
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("E51").Value = "  +1.26%  "

# Numeric-looking text cells in column D: force text format so Excel
# keeps the exact string (trailing zeros, etc.) instead of parsing it as
# a number, then restore the default style so no stray style index is
# left attached to the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "264.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3250"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06804"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7828"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07784"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.033"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007987"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.644"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.477"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.016"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.685"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "111.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.118"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08720"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04867"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7201"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.873"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.120"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.259"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4882"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9017"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.970"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.683"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4203"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.019"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1238"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.05"
$ws.Range("D51").Style = "Normal"
